$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.809.91"
$ws.Range("E2").Value = "  +1.51%  "

# Row 3
$ws.Range("D3").Value = "3.472.77"
$ws.Range("E3").Value = "  +1.77%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").Value = "'414.79"
$ws.Range("E5").Value = "  +1.18%  "

# Row 6
$ws.Range("D6").Value = "'130.66"
$ws.Range("E6").Value = "  +1.33%  "

# Row 7
$ws.Range("D7").Value = "'0.628"
$ws.Range("E7").Value = "  -1.44%  "

# Row 8
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
$ws.Range("D9").Value = "'0.728"
$ws.Range("E9").Value = "  -1.30%  "

# Row 10
$ws.Range("D10").Value = "'0.151"
$ws.Range("E10").Value = "  +5.44%  "

# Row 11
$ws.Range("D11").Value = "'42.64"
$ws.Range("E11").Value = "  -2.23%  "

# Row 12
$ws.Range("D12").Value = "'9.75"
$ws.Range("E12").Value = "  +3.54%  "

# Row 13
$ws.Range("D13").Value = "'0.0000221"
$ws.Range("E13").Value = "  -2.54%  "

# Row 14
$ws.Range("D14").Value = "4.024.03"
$ws.Range("E14").Value = "  +1.97%  "

# Row 15
$ws.Range("D15").Value = "'0.140"
$ws.Range("E15").Value = "  -0.24%  "

# Row 16
$ws.Range("D16").Value = "'20.49"
$ws.Range("E16").Value = "  -4.24%  "

# Row 17
$ws.Range("D17").Value = "3.477.49"
$ws.Range("E17").Value = "  +1.66%  "

# Row 18
$ws.Range("D18").Value = "'12.63"
$ws.Range("E18").Value = "  +0.37%  "

# Row 19
$ws.Range("E19").Value = "  +0.16%  "

# Row 20
$ws.Range("D20").Value = "62.790.23"
$ws.Range("E20").Value = "  +1.54%  "

# Row 21
$ws.Range("D21").Value = "'472.07"
$ws.Range("E21").Value = "  +4.32%  "

# Row 22
$ws.Range("D22").Value = "'90.68"
$ws.Range("E22").Value = "  -1.34%  "

# Row 23
$ws.Range("D23").Value = "'3.29"
$ws.Range("E23").Value = "  +2.73%  "

# Row 24
$ws.Range("D24").Value = "'13.18"
$ws.Range("E24").Value = "  -0.45%  "

# Row 25
$ws.Range("D25").Value = "'10.54"
$ws.Range("E25").Value = "  +13.37%  "

# Row 26
$ws.Range("D26").Value = "'3.32"
$ws.Range("E26").Value = "  +0.73%  "

# Row 27
$ws.Range("D27").Value = "'33.45"
$ws.Range("E27").Value = "  +0.04%  "

# Row 28
$ws.Range("E28").Value = "  +0.15%  "

# Row 29
$ws.Range("D29").Value = "'7.56"
$ws.Range("E29").Value = "  -1.52%  "

# Row 31
$ws.Range("E31").Value = "  -3.49%  "

# Row 33
$ws.Range("D33").Value = "'0.112"
$ws.Range("E33").Value = "  -1.63%  "

# Row 34
$ws.Range("D34").Value = "'40.99"
$ws.Range("E34").Value = "  -3.33%  "

# Row 35
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.01%  "

# Row 36
$ws.Range("D36").Value = "'58.34"
$ws.Range("E36").Value = "  +8.14%  "

# Row 37
$ws.Range("D37").Value = "'0.0488"
$ws.Range("E37").Value = "  -3.33%  "

# Row 38
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  +0.11%  "

# Row 39
$ws.Range("D39").Value = "'3.05"
$ws.Range("E39").Value = "  +3.27%  "

# Row 40
$ws.Range("D40").Value = "'2.80"
$ws.Range("E40").Value = "  +7.83%  "

# Row 41
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'148.75"
$ws.Range("E41").Value = "  +3.22%  "

# Row 42
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.135"
$ws.Range("E42").Value = "  -1.67%  "

# Row 43
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").Value = "'0.321"
$ws.Range("E43").Value = "  +0.00%  "

# Row 44
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "'4.44"
$ws.Range("E44").Value = "  +1.12%  "

# Row 45
$ws.Range("D45").Value = "'3.33"
$ws.Range("E45").Value = "  -2.10%  "

# Row 46
$ws.Range("E46").Value = "  +2.24%  "

# Row 47
$ws.Range("D47").Value = "0.0₃0569"
$ws.Range("E47").Value = "  +29.00%  "

# Row 48
$ws.Range("D48").Value = "'2.38"
$ws.Range("E48").Value = "  +9.58%  "

# Row 49
$ws.Range("D49").Value = "'16.42"
$ws.Range("E49").Value = "  -1.83%  "

# Row 50
$ws.Range("D50").Value = "'22.15"
$ws.Range("E50").Value = "  -2.22%  "

# Row 51
$ws.Range("D51").Value = "'0.141"
$ws.Range("E51").Value = "  -7.34%  "
